# Pat Cummins sheet: keep only the header row plus the "Oct 16 2020" (vs
# Mumbai Indians) and "Oct 7 2020" (vs Chennai Super Kings) match rows,
# which end up re-numbered as rows 2 and 3. All the other match rows are
# removed entirely.
#
# Original row -> new row mapping:
#   row 1 (header)                       -> row 1 (unchanged)
#   row 4 (Oct 16 2020 / Mumbai Indians)  -> row 2
#   row 7 (Oct 7 2020 / Chennai Super K.) -> row 3
#   rows 2,3,5,6,8,9                      -> deleted

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-up so earlier row numbers stay valid as later deletes happen.
$ws.Rows("8:9").Delete()
$ws.Rows("5:6").Delete()
$ws.Rows("2:3").Delete()
